$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 50: header row (same look as row 1), D = Test Error ----
$ws.Range("A50").Value = "Date"
$ws.Range("B50").Value = "Method"
$ws.Range("C50").Value = "Groups"
$ws.Range("D50").Value = "Test Error"
$ws.Range("E50").Value = "Notes"
$hdrRange = $ws.Range("A50:E50")
$hdrRange.Font.Size = 16
$hdrRange.Font.Bold = $true
$hdrRange.Font.ThemeColor = 2
$hdrRange.Interior.Color = 0xB7 + 0xAD*256 + 0xF5*65536
$hdrRange.HorizontalAlignment = -4108

# ---- Regression Tree block: rows 51-58 (cols A-D) ----
$ws.Range("A51").Font.Size = 12
$ws.Range("A51").Font.Bold = $true
$ws.Range("A51").NumberFormat = "m/d/yy"
$ws.Range("A51").HorizontalAlignment = -4108
$ws.Range("A52:A58").Font.Size = 12

$ws.Range("B51").Value = "Regression "
$ws.Range("B52").Value = "Tree"
$regTree = $ws.Range("B51:B52")
$regTree.Font.Size = 12
$regTree.Font.Bold = $true
$regTree.Interior.Color = 0xE2 + 0x14*256 + 0x98*65536
$regTree.HorizontalAlignment = -4108

$ws.Range("B53:B58").Font.Size = 12
$ws.Range("B53:B58").HorizontalAlignment = -4108

$ws.Range("C51").Value = "Top 5"
$ws.Range("C52").Value = "Top 10"
$ws.Range("C53").Value = "Top 15"
$ws.Range("C54").Value = "Top 20"
$ws.Range("C55").Value = "Top 25"
$ws.Range("C56").Value = "Top 30"
$ws.Range("C57").Value = "Top 35"
$ws.Range("C58").Value = "Top 37"
$ws.Range("C51:C58").Font.Size = 12
$ws.Range("C51:C58").HorizontalAlignment = -4108

$ws.Range("D51:D58").Font.Size = 12
$ws.Range("D51:D58").HorizontalAlignment = -4108

# ---- Bagging block: rows 59-66 (cols B,C) ----
$ws.Range("B59").Value = "Bagging "
$ws.Range("B59").Font.Size = 12
$ws.Range("B59").Font.Bold = $true
$ws.Range("B59").Interior.Color = 0xD2 + 0x1C*256 + 0x1C*65536
$ws.Range("B59").HorizontalAlignment = -4108

$ws.Range("B60").Font.Size = 12
$ws.Range("B60").Font.Bold = $true
$ws.Range("B60").HorizontalAlignment = -4108

$ws.Range("B61:B66").Font.Size = 12
$ws.Range("B61:B66").HorizontalAlignment = -4108

$ws.Range("C59").Value = "Top 5"
$ws.Range("C60").Value = "Top 10"
$ws.Range("C61").Value = "Top 15"
$ws.Range("C62").Value = "Top 20"
$ws.Range("C63").Value = "Top 25"
$ws.Range("C64").Value = "Top 30"
$ws.Range("C65").Value = "Top 35"
$ws.Range("C66").Value = "Top 37"
$ws.Range("C59:C66").Font.Size = 12
$ws.Range("C59:C66").HorizontalAlignment = -4108

# ---- Random Forest block: rows 67-74 (cols B,C) ----
$ws.Range("B67").Value = "Random "
$ws.Range("B68").Value = "Forest"
$randForest = $ws.Range("B67:B68")
$randForest.Font.Size = 12
$randForest.Font.Bold = $true
$randForest.Interior.Color = 0xE2 + 0x00*256 + 0x00*65536
$randForest.HorizontalAlignment = -4108

$ws.Range("B69:B74").Font.Size = 12
$ws.Range("B69:B74").HorizontalAlignment = -4108

$ws.Range("C67").Value = "Top 5"
$ws.Range("C68").Value = "Top 10"
$ws.Range("C69").Value = "Top 15"
$ws.Range("C70").Value = "Top 20"
$ws.Range("C71").Value = "Top 25"
$ws.Range("C72").Value = "Top 30"
$ws.Range("C73").Value = "Top 35"
$ws.Range("C74").Value = "Top 37"
$ws.Range("C67:C74").Font.Size = 12
$ws.Range("C67:C74").HorizontalAlignment = -4108

# ---- Boosting block: rows 75-82 (cols B,C) ----
$ws.Range("B75").Value = "Boosting "
$ws.Range("B75").Font.Size = 12
$ws.Range("B75").Font.Bold = $true
$ws.Range("B75").Interior.ThemeColor = 8
$ws.Range("B75").HorizontalAlignment = -4108

$ws.Range("B76").Font.Size = 12
$ws.Range("B76").Font.Bold = $true
$ws.Range("B76").HorizontalAlignment = -4108

$ws.Range("B77:B82").Font.Size = 12
$ws.Range("B77:B82").HorizontalAlignment = -4108

$ws.Range("C75").Value = "Top 5"
$ws.Range("C76").Value = "Top 10"
$ws.Range("C77").Value = "Top 15"
$ws.Range("C78").Value = "Top 20"
$ws.Range("C79").Value = "Top 25"
$ws.Range("C80").Value = "Top 30"
$ws.Range("C81").Value = "Top 35"
$ws.Range("C82").Value = "Top 37"
$ws.Range("C75:C82").Font.Size = 12
$ws.Range("C75:C82").HorizontalAlignment = -4108

# ---- Update view state to match target ----
$excel.ActiveWindow.ScrollRow = 40
$ws.Range("E68").Select()
